# Modified 05/08/2024 03:42:59 PM IST
# Resize the data sheet's columns. The underlying OOXML stores widths that
# are 5/6 (0.8333...) wider than the COM "ColumnWidth" value, so the
# ColumnWidth values below are chosen such that after Excel's internal
# padding conversion the saved <col width="..."> values end up exactly:
# 12, 9, 7, 8, 8, 7, 5, 13, 6, 50 for columns A..J respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth  = 11.166666666666666   # -> stored width 12
$ws.Columns.Item(2).ColumnWidth  = 8.166666666666666    # -> stored width 9
$ws.Columns.Item(3).ColumnWidth  = 6.166666666666667    # -> stored width 7
$ws.Columns.Item(4).ColumnWidth  = 7.166666666666667    # -> stored width 8
$ws.Columns.Item(5).ColumnWidth  = 7.166666666666667    # -> stored width 8
$ws.Columns.Item(6).ColumnWidth  = 6.166666666666667    # -> stored width 7
$ws.Columns.Item(7).ColumnWidth  = 4.166666666666667    # -> stored width 5
$ws.Columns.Item(8).ColumnWidth  = 12.166666666666666   # -> stored width 13
$ws.Columns.Item(9).ColumnWidth  = 5.166666666666667    # -> stored width 6
$ws.Columns.Item(10).ColumnWidth = 49.166666666666664   # -> stored width 50
